$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 43
$ws.Range("A43").Value = 4
$ws.Range("B43").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C43").Value = "Los Lagos"
$ws.Range("D43").Value = 44939
$ws.Range("D43").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E43").Value = 10
$ws.Range("F43").Value = "Fruta"
$ws.Range("G43").Value = 100101
$ws.Range("H43").Value = "Berries"
$ws.Range("I43").Value = 100101001
$ws.Range("J43").Value = "Arándano (blue)"
$ws.Range("K43").Value = "Sin especificar"
$ws.Range("L43").Value = "Primera"
$ws.Range("M43").Value = 200
$ws.Range("N43").Value = 2000
$ws.Range("O43").Value = 2000
$ws.Range("P43").Value = 2000
$ws.Range("Q43").Value = "`$/bandeja 2 kilos"
$ws.Range("R43").Value = "Provincia de Curicó"
$ws.Range("S43").Value = 1000
$ws.Range("T43").Value = 2

# New row 44
$ws.Range("A44").Value = 4
$ws.Range("B44").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C44").Value = "Los Lagos"
$ws.Range("D44").Value = 44939
$ws.Range("D44").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E44").Value = 10
$ws.Range("F44").Value = "Fruta"
$ws.Range("G44").Value = 100101
$ws.Range("H44").Value = "Berries"
$ws.Range("I44").Value = 100101001
$ws.Range("J44").Value = "Arándano (blue)"
$ws.Range("K44").Value = "Sin especificar"
$ws.Range("L44").Value = "Segunda"
$ws.Range("M44").Value = 200
$ws.Range("N44").Value = 2200
$ws.Range("O44").Value = 2200
$ws.Range("P44").Value = 2200
$ws.Range("Q44").Value = "`$/bandeja 2 kilos"
$ws.Range("R44").Value = "Provincia de Curicó"
$ws.Range("S44").Value = 1100
$ws.Range("T44").Value = 2
